$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record at row 264 (pushes the existing rows 264-285 down to 265-286,
# growing the used range from A1:R285 to A1:R286).
$ws.Rows.Item(264).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A264").Value = 11
$ws.Range("B264").Value = "Vega Monumental Concepción"
$ws.Range("C264").Value = "Bíobío"
$ws.Range("D264").Value = 45223
$ws.Range("E264").Value = 8
$ws.Range("F264").Value = 100112032
$ws.Range("G264").Value = "Zapallo italiano"
$ws.Range("H264").Value = "Sin especificar"
$ws.Range("I264").Value = "Primera"
$ws.Range("J264").Value = 150
$ws.Range("K264").Value = 18000
$ws.Range("L264").Value = 18000
$ws.Range("M264").Value = 18000
$ws.Range("N264").Value = "`$/caja 50 unidades"
$ws.Range("O264").Value = "Región de Arica y Parinacota"
$ws.Range("P264").Value = 360
$ws.Range("Q264").Value = 50
$ws.Range("R264").Value = "Hortaliza"
